$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated LR-pair values ("Natmi following Dr Hou advice"):
# shared strings gain a new "FAPs" cluster entry and every sending/target cluster
# combination (ECs, FAPs, M2, sCs) x (ECs, FAPs, M2, sCs) is now represented,
# expanding the table from 12 to 16 data rows (rows 2-17).
$data = @(
  @("ECs","Pdgfb","Pdgfra","ECs",3,1,[double]"53.435824",[double]"160.307472",[double]"0.8055519732580867",[double]"0.8055519732580868",3,1,[double]"3.535386",[double]"10.606158",[double]"0.01988747852527457",[double]"0.01988747852527457",[double]"188.916264068064",[double]"1700.246376612576",[double]"0.01602039756916276",[double]"0.01602039756916275"),
  @("ECs","Pdgfb","Pdgfra","FAPs",3,1,[double]"53.435824",[double]"160.307472",[double]"0.8055519732580867",[double]"0.8055519732580868",3,1,[double]"173.8189136666666",[double]"521.456741",[double]"0.9777772251268709",[double]"0.9777772251268707",[double]"9288.156878563193",[double]"83593.41190706876",[double]"0.7876503731077673",[double]"0.7876503731077673"),
  @("ECs","Pdgfb","Pdgfra","M2",3,1,[double]"53.435824",[double]"160.307472",[double]"0.8055519732580867",[double]"0.8055519732580868",3,1,[double]"0.06908833333333333",[double]"0.207265",[double]"0.000388640093475982",[double]"0.0003886400934759819",[double]"3.691792020453334",[double]"33.22612818408",[double]"0.0003130697941867846",[double]"0.0003130697941867845"),
  @("ECs","Pdgfb","Pdgfra","sCs",3,1,[double]"53.435824",[double]"160.307472",[double]"0.8055519732580867",[double]"0.8055519732580868",3,1,[double]"0.346056",[double]"1.038168",[double]"0.001946656254378565",[double]"0.001946656254378564",[double]"18.491787510144",[double]"166.426087591296",[double]"0.001568132786969849",[double]"0.001568132786969849"),
  @("FAPs","Pdgfb","Pdgfra","ECs",1,[double]"0.3333333333333333",[double]"0.03171066666666666",[double]"0.095132",[double]"0.0004780424103998614",[double]"0.0004780424103998615",3,1,[double]"3.535386",[double]"10.606158",[double]"0.01988747852527457",[double]"0.01988747852527457",[double]"0.112109446984",[double]"1.008985022856",[double]"9.507058170997738e-06",[double]"9.507058170997737e-06"),
  @("FAPs","Pdgfb","Pdgfra","FAPs",1,[double]"0.3333333333333333",[double]"0.03171066666666666",[double]"0.095132",[double]"0.0004780424103998614",[double]"0.0004780424103998615",3,1,[double]"173.8189136666666",[double]"521.456741",[double]"0.9777772251268709",[double]"0.9777772251268707",[double]"5.511913631645776",[double]"49.60722268481199",[double]"0.0004674189815337373",[double]"0.0004674189815337373"),
  @("FAPs","Pdgfb","Pdgfra","M2",1,[double]"0.3333333333333333",[double]"0.03171066666666666",[double]"0.095132",[double]"0.0004780424103998614",[double]"0.0004780424103998615",3,1,[double]"0.06908833333333333",[double]"0.207265",[double]"0.000388640093475982",[double]"0.0003886400934759819",[double]"0.002190837108888889",[double]"0.01971753398",[double]"1.857864470632859e-07",[double]"1.857864470632859e-07"),
  @("FAPs","Pdgfb","Pdgfra","sCs",1,[double]"0.3333333333333333",[double]"0.03171066666666666",[double]"0.095132",[double]"0.0004780424103998614",[double]"0.0004780424103998615",3,1,[double]"0.346056",[double]"1.038168",[double]"0.001946656254378565",[double]"0.001946656254378564",[double]"0.010973666464",[double]"0.098762998176",[double]"9.305842480630949e-07",[double]"9.305842480630948e-07"),
  @("M2","Pdgfb","Pdgfra","ECs",3,1,[double]"10.73549533333333",[double]"32.206486",[double]"0.1618389837063174",[double]"0.1618389837063175",3,1,[double]"3.535386",[double]"10.606158",[double]"0.01988747852527457",[double]"0.01988747852527457",[double]"37.954119904532",[double]"341.587079140788",[double]"0.00321856931301165",[double]"0.003218569313011649"),
  @("M2","Pdgfb","Pdgfra","FAPs",3,1,[double]"10.73549533333333",[double]"32.206486",[double]"0.1618389837063174",[double]"0.1618389837063175",3,1,[double]"173.8189136666666",[double]"521.456741",[double]"0.9777772251268709",[double]"0.9777772251268707",[double]"1866.032136513569",[double]"16794.28922862212",[double]"0.158242472405716",[double]"0.1582424724057159"),
  @("M2","Pdgfb","Pdgfra","M2",3,1,[double]"10.73549533333333",[double]"32.206486",[double]"0.1618389837063174",[double]"0.1618389837063175",3,1,[double]"0.06908833333333333",[double]"0.207265",[double]"0.000388640093475982",[double]"0.0003886400934759819",[double]"0.7416974800877778",[double]"6.675277320789999",[double]"6.289711775568113e-05",[double]"6.289711775568113e-05"),
  @("M2","Pdgfb","Pdgfra","sCs",3,1,[double]"10.73549533333333",[double]"32.206486",[double]"0.1618389837063174",[double]"0.1618389837063175",3,1,[double]"0.346056",[double]"1.038168",[double]"0.001946656254378565",[double]"0.001946656254378564",[double]"3.715082573071999",[double]"33.435743157648",[double]"0.0003150448698341735",[double]"0.0003150448698341735"),
  @("sCs","Pdgfb","Pdgfra","ECs",3,1,[double]"2.131391333333333",[double]"6.394174",[double]"0.03213100062519576",[double]"0.03213100062519577",3,1,[double]"3.535386",[double]"10.606158",[double]"0.01988747852527457",[double]"0.01988747852527457",[double]"7.535291080388",[double]"67.817619723492",[double]"0.0006390045849291647",[double]"0.0006390045849291645"),
  @("sCs","Pdgfb","Pdgfra","FAPs",3,1,[double]"2.131391333333333",[double]"6.394174",[double]"0.03213100062519576",[double]"0.03213100062519577",3,1,[double]"173.8189136666666",[double]"521.456741",[double]"0.9777772251268709",[double]"0.9777772251268707",[double]"370.4761261585481",[double]"3334.285135426934",[double]"0.03141696063185367",[double]"0.03141696063185367"),
  @("sCs","Pdgfb","Pdgfra","M2",3,1,[double]"2.131391333333333",[double]"6.394174",[double]"0.03213100062519576",[double]"0.03213100062519577",3,1,[double]"0.06908833333333333",[double]"0.207265",[double]"0.000388640093475982",[double]"0.0003886400934759819",[double]"0.1472542749011111",[double]"1.32528847411",[double]"1.248739508645292e-05",[double]"1.248739508645292e-05"),
  @("sCs","Pdgfb","Pdgfra","sCs",3,1,[double]"2.131391333333333",[double]"6.394174",[double]"0.03213100062519576",[double]"0.03213100062519577",3,1,[double]"0.346056",[double]"1.038168",[double]"0.001946656254378565",[double]"0.001946656254378564",[double]"0.7375807592479999",[double]"6.638226833231999",[double]"6.25480133264789e-05",[double]"6.25480133264789e-05"),
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c - 1]
  }
}
